$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 15307.85069381629
$ws.Range("D3").Value = 833.5794467324417
$ws.Range("E3").Value = 2825.673217592491

$ws.Range("B4").Value = 7241.871470626517
$ws.Range("D4").Value = 442.1657644405556
$ws.Range("E4").Value = 2579.961218619072

$ws.Range("B5").Value = 2271.016
$ws.Range("E5").Value = 22.023

$ws.Range("B6").Value = 9328.085500000008
$ws.Range("D6").Value = 586.004
$ws.Range("E6").Value = 1282.001

$ws.Range("B7").Value = 13779.41100000003
$ws.Range("D7").Value = 810
$ws.Range("E7").Value = 2076.004

$ws.Range("B8").Value = 21117.26400000003
$ws.Range("D8").Value = 1095
$ws.Range("E8").Value = 3380.005

$ws.Range("B9").Value = 48357.26000000002
$ws.Range("D9").Value = 4532.005
$ws.Range("E9").Value = 22752.004

$ws.Range("F10").Value = 8045806324.673006

$ws.Range("G11").Value = 0.760955816886735

$ws.Range("F12").Value = 438129357.203
$ws.Range("G12").Value = 0.05445437530101201

$ws.Range("F13").Value = 1485173843.166
$ws.Range("G13").Value = 0.184589807812253
